$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.780058622360229
$ws.Range("B1").Value = 1.862804532051086
$ws.Range("C1").Value = 2.038267850875854
$ws.Range("D1").Value = 2.899649858474731
$ws.Range("E1").Value = 4.799004554748535
